$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.310.64"
$ws.Range("D3").Value = "2.634.97"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'604.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").Value = "'155.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  +8.84%  "
$ws.Range("E10").Value = "  +5.77%  "
$ws.Range("D11").Value = "'5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "'29.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +23.41%  "
$ws.Range("D15").Value = "3.107.56"
$ws.Range("D16").Value = "65.170.43"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "2.637.74"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'12.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("D19").Value = "'4.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'359.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("D21").Value = "'7.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.93%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'69.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").Value = "'9.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "'1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'8.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").Value = "0.0₃0957"
$ws.Range("E29").Value = "  +13.41%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'2.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.90%  "
$ws.Range("D32").Value = "'525.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "'5.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.38%  "
$ws.Range("D35").Value = "'6.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.54%  "
$ws.Range("D36").Value = "'0.428"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").Value = "'20.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("D39").Value = "'162.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'42.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("D43").Value = "'165.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'4.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").Value = "'0.0618"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.70%  "
$ws.Range("D46").Value = "'23.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("E47").Value = "  +5.92%  "
$ws.Range("E48").Value = "  +7.76%  "
$ws.Range("D49").Value = "'0.655"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.04%  "
$ws.Range("D50").Value = "'0.0982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "'19.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
